# Update the "Ccl21b-Cxcr3" LR-pair sheet with refreshed TPM-derived values.
# The new data adds an "ECs" sending-cluster pair (rows 2-3), shifts the
# previous "FAPs" pair down to rows 4-5, and appends a new "MuSCs" pair
# (rows 6-7), each with recalculated specificity/expression metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 6,20

$arr[0,0] = "ECs"
$arr[0,1] = "Ccl21b"
$arr[0,2] = "Cxcr3"
$arr[0,3] = "FAPs"
$arr[0,4] = 1
$arr[0,5] = 0.3333333333333333
$arr[0,6] = 0.004739
$arr[0,7] = 0.014217
$arr[0,8] = 0.02588570741885795
$arr[0,9] = 0.02588570741885795
$arr[0,10] = 1
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.05194133333333333
$arr[0,13] = 0.155824
$arr[0,14] = 0.02773017886769741
$arr[0,15] = 0.02773017886769741
$arr[0,16] = 0.0002461499786666667
$arr[0,17] = 0.002215349808
$arr[0,18] = 0.0007178152968418128
$arr[0,19] = 0.0007178152968418129
$arr[1,0] = "ECs"
$arr[1,1] = "Ccl21b"
$arr[1,2] = "Cxcr3"
$arr[1,3] = "Resolving-Mac"
$arr[1,4] = 1
$arr[1,5] = 0.3333333333333333
$arr[1,6] = 0.004739
$arr[1,7] = 0.014217
$arr[1,8] = 0.02588570741885795
$arr[1,9] = 0.02588570741885795
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 1.821156333333333
$arr[1,13] = 5.463469
$arr[1,14] = 0.9722698211323025
$arr[1,15] = 0.9722698211323026
$arr[1,16] = 0.008630459863666668
$arr[1,17] = 0.077674138773
$arr[1,18] = 0.02516789212201613
$arr[1,19] = 0.02516789212201614
$arr[2,0] = "FAPs"
$arr[2,1] = "Ccl21b"
$arr[2,2] = "Cxcr3"
$arr[2,3] = "FAPs"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.134289
$arr[2,7] = 0.402867
$arr[2,8] = 0.733523056250478
$arr[2,9] = 0.733523056250478
$arr[2,10] = 1
$arr[2,11] = 0.3333333333333333
$arr[2,12] = 0.05194133333333333
$arr[2,13] = 0.155824
$arr[2,14] = 0.02773017886769741
$arr[2,15] = 0.02773017886769741
$arr[2,16] = 0.006975149712000001
$arr[2,17] = 0.06277634740800001
$arr[2,18] = 0.02034072555340583
$arr[2,19] = 0.02034072555340583
$arr[3,0] = "FAPs"
$arr[3,1] = "Ccl21b"
$arr[3,2] = "Cxcr3"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 0.134289
$arr[3,7] = 0.402867
$arr[3,8] = 0.733523056250478
$arr[3,9] = 0.733523056250478
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 1.821156333333333
$arr[3,13] = 5.463469
$arr[3,14] = 0.9722698211323025
$arr[3,15] = 0.9722698211323026
$arr[3,16] = 0.2445612628470001
$arr[3,17] = 2.201051365623
$arr[3,18] = 0.7131823306970722
$arr[3,19] = 0.7131823306970722
$arr[4,0] = "MuSCs"
$arr[4,1] = "Ccl21b"
$arr[4,2] = "Cxcr3"
$arr[4,3] = "FAPs"
$arr[4,4] = 2
$arr[4,5] = 0.6666666666666666
$arr[4,6] = 0.04404599999999999
$arr[4,7] = 0.132138
$arr[4,8] = 0.2405912363306641
$arr[4,9] = 0.2405912363306641
$arr[4,10] = 1
$arr[4,11] = 0.3333333333333333
$arr[4,12] = 0.05194133333333333
$arr[4,13] = 0.155824
$arr[4,14] = 0.02773017886769741
$arr[4,15] = 0.02773017886769741
$arr[4,16] = 0.002287807968
$arr[4,17] = 0.02059027171199999
$arr[4,18] = 0.006671638017449775
$arr[4,19] = 0.006671638017449775
$arr[5,0] = "MuSCs"
$arr[5,1] = "Ccl21b"
$arr[5,2] = "Cxcr3"
$arr[5,3] = "Resolving-Mac"
$arr[5,4] = 2
$arr[5,5] = 0.6666666666666666
$arr[5,6] = 0.04404599999999999
$arr[5,7] = 0.132138
$arr[5,8] = 0.2405912363306641
$arr[5,9] = 0.2405912363306641
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 1.821156333333333
$arr[5,13] = 5.463469
$arr[5,14] = 0.9722698211323025
$arr[5,15] = 0.9722698211323026
$arr[5,16] = 0.08021465185799999
$arr[5,17] = 0.7219318667219998
$arr[5,18] = 0.2339195983132143
$arr[5,19] = 0.2339195983132143

$ws.Range("A2:T7").Value = $arr
